$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7287194209349384
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 9.576116808119359

$ws.Range("B3").Value = 0.7287194209349384
$ws.Range("C3").Value = 9.226618575922256
$ws.Range("D3").Value = 0.7127328510149897
$ws.Range("E3").Value = 6.48142807727062
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 17.1494989251428

$ws.Range("B4").Value = 3.182878228561681
$ws.Range("C4").Value = 1.65323645889881
$ws.Range("D4").Value = 0.1529057820181812
$ws.Range("E4").Value = 0.4998867070740569
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 5.488907176552729

$ws.Range("B5").Value = 1.505614041169197
$ws.Range("C5").Value = 1.65323645889881
$ws.Range("D5").Value = 0.7127328510149897
$ws.Range("E5").Value = 0.4998867070740569
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.371470058157054

$ws.Range("B6").Value = 1.505614041169197
$ws.Range("C6").Value = 1.65323645889881
$ws.Range("D6").Value = 0.1529057820181812
$ws.Range("E6").Value = 0.4998867070740569
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 3.811642989160245

$ws.Range("B7").Value = 0.1554434735375247
$ws.Range("C7").Value = 9.226618575922256
$ws.Range("D7").Value = 157.8057217802531
$ws.Range("E7").Value = 6.48142807727062
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 173.6692119069835

$ws.Range("B8").Value = 3.182878228561681
$ws.Range("C8").Value = 1766.335244827366
$ws.Range("D8").Value = 157.8057217802531
$ws.Range("E8").Value = 6.48142807727062
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1933.805272913452

$ws.Range("B9").Value = 0.1554434735375247
$ws.Range("C9").Value = 0.3375848360084654
$ws.Range("D9").Value = 16.98373111632243
$ws.Range("E9").Value = 6.48142807727062
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 23.95818750313904

$ws.Range("B10").Value = 0.06328177979961902
$ws.Range("C10").Value = 9.226618575922256
$ws.Range("D10").Value = 0.7127328510149897
$ws.Range("E10").Value = 6.48142807727062
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 16.48406128400748

$ws.Range("B11").Value = 0.3464964993005633
$ws.Range("C11").Value = 1.65323645889881
$ws.Range("D11").Value = 16.98373111632243
$ws.Range("E11").Value = 6.48142807727062
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 25.46489215179242
